$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R data for year 2021
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

$ws.Range("M5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 18.953297329007047

# Update the selection to match the new active cell
$ws.Range("Q8").Select()
